$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.008.14'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '2.356.32'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  -0.09%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '499.24'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -1.90%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.35'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -3.42%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').Value = '2.355.32'
$ws.Range('E9').Value = '  -3.79%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '2.776.48'
$ws.Range('E14').Value = '  -3.74%  '
$ws.Range('D15').Value = '55.957.09'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('E16').Value = '  -2.65%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000131'
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.420.06'
$ws.Range('E18').Value = '  -2.09%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.95'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('E20').Value = '  -2.50%  '
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '305.23'
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  -2.83%  '
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('E23').Value = '  +0.04%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.05'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -0.15%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.995'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('E26').Value = '  -3.61%  '
$ws.Range('E27').Value = '  -6.50%  '
$ws.Range('E28').Value = '  -4.98%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.25'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('E30').Value = '  -3.74%  '
$style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.63'
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('E32').Value = '  +0.07%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.71'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  -7.32%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.07'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -5.26%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.53'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  -2.63%  '
$style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.16'
$ws.Range('D37').Style = $style
$ws.Range('E37').Value = '  -6.07%  '
$ws.Range('E38').Value = '  -3.05%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.99'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -1.94%  '
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('E41').Value = '  -6.16%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '128.77'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('E43').Value = '  -1.87%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.65'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  -4.84%  '
$ws.Range('E45').Value = '  -2.22%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0899'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  -1.97%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '238.06'
$ws.Range('D47').Style = $style
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('E49').Value = '  -3.65%  '
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('E51').Value = '  -0.68%  '
